$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - force text format to preserve exact string representation
# (values like "0.140", "599.91" etc. would otherwise be auto-converted to numbers by Excel)
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.730.48'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.711.24'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '599.91'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '163.11'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.544'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.710.54'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.140'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.45'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.208.16'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '68.693.74'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.692.97'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.86'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '365.26'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '74.01'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.92'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.840.30'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '597.11'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.999'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.45'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.88'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '160.59'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.380'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.44'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0₆0317'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '158.11'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '22.13'

# Column E (Volume(1h)) updates
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('E6').Value = '  +3.74%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E10').Value = '  -0.74%  '
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('E13').Value = '  +2.81%  '
$ws.Range('E14').Value = '  +1.43%  '
$ws.Range('E15').Value = '  +2.49%  '
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('E18').Value = '  +1.40%  '
$ws.Range('E19').Value = '  +4.14%  '
$ws.Range('E20').Value = '  +4.73%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('E23').Value = '  +2.50%  '
$ws.Range('E24').Value = '  +2.99%  '
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('E28').Value = '  +2.27%  '
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('E30').Value = '  +6.68%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('E32').Value = '  +3.05%  '
$ws.Range('E33').Value = '  +3.24%  '
$ws.Range('E34').Value = '  +4.91%  '
$ws.Range('E35').Value = '  +2.97%  '
$ws.Range('E36').Value = '  +5.64%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('E40').Value = '  +2.24%  '
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('E42').Value = '  +1.85%  '
$ws.Range('E43').Value = '  +3.23%  '
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('E45').Value = '  -5.24%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('E47').Value = '  -0.51%  '
$ws.Range('E48').Value = '  +5.33%  '
$ws.Range('E49').Value = '  +6.06%  '
$ws.Range('E50').Value = '  +7.28%  '
$ws.Range('E51').Value = '  -0.21%  '
